# "Update on monitoring report" -- adds a "strata" column (O) to the
# sampling frame, mirroring the existing "Location" column (N) but using
# underscore-joined values (e.g. "East_Jerusalem") instead of spaced text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("O1").Value = "strata"

# Fill every data row (2-35) with the underscored strata value that
# corresponds to column N's "East Jerusalem" text.
$lastRow = 35
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = "East_Jerusalem"
}

# Reflect the author's final on-screen selection.
$ws.Range("N35").Select()
